$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "#Trump announced that he would fire #WilliamBarr.`n.`n.`nHe reminds me of children ...😑`n#TrumpIsACriminal… https://t.co/1B51buQrc6"
$ws.Range("C2").Value = 1338856399380623000
$ws.Range("D2").Value = 128
$ws.Range("E2").Value = 44180.61128472222
$ws.Range("F2").Value = "en"
$ws.Range("G2").Value = 1273713204242395000
$ws.Range("H2").Value = "Ashleysullivan_"
$ws.Range("I2").Value = 596
$ws.Range("J2").Value = "Chicago"
$ws.Range("K2").Value = "#Trump announced that he would fire #WilliamBarr...He reminds me of children ...😑#TrumpIsACriminal… https://t.co/1B51buQrc6 "
$ws.Range("L2").Value = "#Trump announced that he would fire #WilliamBarr.`n.`n.`nHe reminds me of children ...😑`n#TrumpIsACriminal… "
$ws.Range("M2").Value = "#Trump announced that he would fire #WilliamBarr...He reminds me of children ...😑#TrumpIsACriminal… "

# Row 3
$ws.Range("B3").Value = "RT @ManDessins: #COVID19 #USA #Trump #Dessin @Midilibre https://t.co/hv8DQwBzIc"
$ws.Range("C3").Value = 1338856388823400000
$ws.Range("D3").Value = 79
$ws.Range("E3").Value = 44180.61125
$ws.Range("F3").Value = "und"
$ws.Range("G3").Value = 604704562
$ws.Range("H3").Value = "bonsoirmichel"
$ws.Range("I3").Value = 13262
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "RT @ManDessins: # COVID19 #USA #Trump #Dessin @Midilibre https://t.co/hv8DQwBzIc "
$ws.Range("L3").Value = " #COVID19 #USA #Trump #Dessin  "
$ws.Range("M3").Value = " # COVID19 #USA #Trump #Dessin   "

# Row 4
$ws.Range("B4").Value = "#Trump `nthis is the end`nJim Morrison"
$ws.Range("C4").Value = 1338856381877604000
$ws.Range("D4").Value = 37
$ws.Range("E4").Value = 44180.61122685186
$ws.Range("F4").Value = "en"
$ws.Range("G4").Value = 16604304
$ws.Range("H4").Value = "muehlenwind"
$ws.Range("I4").Value = 2340
$ws.Range("J4").Value = "Sandstrand, Deutschland"
$ws.Range("K4").Value = "#Trump this is the endJim Morrison "
$ws.Range("L4").Value = "#Trump `nthis is the end`nJim Morrison"
$ws.Range("M4").Value = "#Trump this is the endJim Morrison "

# Row 5
$ws.Range("B5").Value = "@NotHoodlum Why is everyone so shocked? It’s the one promise #Trump kept - the entitled get more entitlement.… https://t.co/I95RVALzIQ"
$ws.Range("C5").Value = 1338856364635009000
$ws.Range("D5").Value = 134
$ws.Range("E5").Value = 44180.61118055556
$ws.Range("F5").Value = "en"
$ws.Range("G5").Value = 757529364
$ws.Range("H5").Value = "MillieMinet"
$ws.Range("I5").Value = 615
$ws.Range("J5").Value = "NPT"
$ws.Range("K5").Value = "@NotHoodlum Why is everyone so shocked? It’s the one promise #Trump kept - the entitled get more entitlement.… https://t.co/I95RVALzIQ "
$ws.Range("L5").Value = " Why is everyone so shocked? It’s the one promise #Trump kept - the entitled get more entitlement.… "
$ws.Range("M5").Value = " Why is everyone so shocked? It’s the one promise #Trump kept - the entitled get more entitlement.…  "

# Row 6
$ws.Range("B6").Value = "President Donald #Trump stands among Army cadets as he attends the annual Army-Navy collegiate football game in… https://t.co/P7nglPVrEb"
$ws.Range("C6").Value = 1338856361669616000
$ws.Range("D6").Value = 136
$ws.Range("E6").Value = 44180.61118055556
$ws.Range("F6").Value = "en"
$ws.Range("G6").Value = 842343300115927000
$ws.Range("H6").Value = "pow_photos"
$ws.Range("I6").Value = 806
$ws.Range("J6").Value = "Lebanon"
$ws.Range("K6").Value = "President Donald #Trump stands among Army cadets as he attends the annual Army-Navy collegiate football game in… https://t.co/P7nglPVrEb "
$ws.Range("L6").Value = "President Donald #Trump stands among Army cadets as he attends the annual Army-Navy collegiate football game in… "
$ws.Range("M6").Value = "President Donald #Trump stands among Army cadets as he attends the annual Army-Navy collegiate football game in…  "

# Row 7
$ws.Range("B7").Value = "The electoral college has spoken. YOU LOST AGAIN. Perhaps #Trump is bad for America. `n#ShutHimOut… https://t.co/1w7CuvK9Ri"
$ws.Range("C7").Value = 1338856349032178000
$ws.Range("D7").Value = 123
$ws.Range("E7").Value = 44180.61114583333
$ws.Range("F7").Value = "en"
$ws.Range("G7").Value = 267129243
$ws.Range("H7").Value = "TarikuBogale"
$ws.Range("I7").Value = 5466
$ws.Range("J7").Value = "New York"
$ws.Range("K7").Value = "The electoral college has spoken. YOU LOST AGAIN. Perhaps #Trump is bad for America. #ShutHimOut… https://t.co/1w7CuvK9Ri "
$ws.Range("L7").Value = "The electoral college has spoken. YOU LOST AGAIN. Perhaps #Trump is bad for America. `n#ShutHimOut… "
$ws.Range("M7").Value = "The electoral college has spoken. YOU LOST AGAIN. Perhaps #Trump is bad for America. #ShutHimOut…  "

# Row 8
$ws.Range("B8").Value = "RT @AlexdGtze: @washingtonpost Truth hurts`n#BidenVaccine #Trump https://t.co/1e8oIOxc5y"
$ws.Range("C8").Value = 1338856347044033000
$ws.Range("D8").Value = 87
$ws.Range("E8").Value = 44180.61113425926
$ws.Range("F8").Value = "en"
$ws.Range("G8").Value = 1259356052635542000
$ws.Range("H8").Value = "JoeKing84576078"
$ws.Range("I8").Value = 511
$ws.Range("J8").Value = "Texas, USA"
$ws.Range("K8").Value = "RT @AlexdGtze: @washingtonpost Truth hurts#BidenVaccine #Trump https://t.co/1e8oIOxc5y "
$ws.Range("L8").Value = "  Truth hurts`n#BidenVaccine #Trump "
$ws.Range("M8").Value = "  Truth hurts#BidenVaccine #Trump  "

# Row 9
$ws.Range("B9").Value = "@mkraju Strongly disagree with vacs for #Trump and #Pence. They have been advocating for super-spreading since the… https://t.co/Jm35wvOwH1"
$ws.Range("C9").Value = 1338856345882096000
$ws.Range("D9").Value = 139
$ws.Range("E9").Value = 44180.61113425926
$ws.Range("F9").Value = "en"
$ws.Range("G9").Value = 1229178950347149000
$ws.Range("H9").Value = "LAShake2"
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = "@mkraju Strongly disagree with vacs for #Trump and #Pence. They have been advocating for super-spreading since the… https://t.co/Jm35wvOwH1 "
$ws.Range("L9").Value = " Strongly disagree with vacs for #Trump and #Pence. They have been advocating for super-spreading since the… "
$ws.Range("M9").Value = " Strongly disagree with vacs for #Trump and #Pence. They have been advocating for super-spreading since the…  "

# Row 10
$ws.Range("B10").Value = "@JoeBiden We all know you cheated. The sad thing is, YOU probably don’t know you cheated. The audacity to write thi… https://t.co/wXz7eT8rMF"
$ws.Range("C10").Value = 1338856318673809000
$ws.Range("D10").Value = 140
$ws.Range("E10").Value = 44180.61105324074
$ws.Range("G10").Value = 1298013110691758000
$ws.Range("H10").Value = "JiruJoshua"
$ws.Range("I10").Value = 17
$ws.Range("J10").Value = "Madison, WI"
$ws.Range("K10").Value = "@JoeBiden We all know you cheated. The sad thing is, YOU probably don’t know you cheated. The audacity to write thi… https://t.co/wXz7eT8rMF "
$ws.Range("L10").Value = " We all know you cheated. The sad thing is, YOU probably don’t know you cheated. The audacity to write thi… "
$ws.Range("M10").Value = " We all know you cheated. The sad thing is, YOU probably don’t know you cheated. The audacity to write thi…  "

# Row 11
$ws.Range("B11").Value = "RT @JUANdeITALIA: El presidente MR #TRUMP sea como sea EL en su conducta imperfecta ...es un patriota además de persona creyente en DIOS, t…"
$ws.Range("C11").Value = 1338856311044321000
$ws.Range("D11").Value = 140
$ws.Range("E11").Value = 44180.61103009259
$ws.Range("F11").Value = "es"
$ws.Range("G11").Value = 353243516
$ws.Range("H11").Value = "no_me_la_calo"
$ws.Range("I11").Value = 1582
$ws.Range("K11").Value = "RT @JUANdeITALIA: The president MR #TRUMP whatever HE may be in his imperfect behavior ... he is a patriot as well as a person who believes in GOD, t ... "
$ws.Range("L11").Value = " El presidente MR #TRUMP sea como sea EL en su conducta imperfecta ...es un patriota además de persona creyente en DIOS, t…"
$ws.Range("M11").Value = " The president MR #TRUMP whatever HE may be in his imperfect behavior ... he is a patriot as well as a person who believes in GOD, t ... "
